$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 25221.834
$ws.Range("I6").Value = 33401.11
$ws.Range("J6").Value = 684
$ws.Range("K6").Value = 100203.33
$ws.Range("L6").Value = 2052
$ws.Range("M6").Value = -100091.33
$ws.Range("N6").Value = -2276
$ws.Range("H33").Value = 234.11111
$ws.Range("I33").Value = 253.21428
$ws.Range("K33").Value = 253.21428
$ws.Range("M33").Value = -24.21428
$ws.Range("H76").Value = 33345998
$ws.Range("J76").Value = 7999
$ws.Range("L76").Value = 7999
$ws.Range("N76").Value = -8629
$ws.Range("H79").Value = 33345998
$ws.Range("J79").Value = 7999
$ws.Range("L79").Value = 7999
$ws.Range("N79").Value = -10183

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H69").Value = 354998.94
$ws.Range("J69").Value = 354998.94
$ws.Range("L69").Value = 354998.94
$ws.Range("N69").Value = -356496.94
$ws.Range("H72").Value = 354998.94
$ws.Range("J72").Value = 354998.94
$ws.Range("L72").Value = 1064996.82
$ws.Range("N72").Value = -1072484.82
$ws.Range("H110").Value = 41668290
$ws.Range("I110").Value = 83333620
$ws.Range("J110").Value = 2963.6667
$ws.Range("K110").Value = 83333620
$ws.Range("L110").Value = 2963.6667
$ws.Range("M110").Value = -83331575
$ws.Range("N110").Value = -7053.6667
$ws.Range("H122").Value = 1857.8572
$ws.Range("I122").Value = 1503.5
$ws.Range("J122").Value = 1999.6
$ws.Range("K122").Value = 4510.5
$ws.Range("L122").Value = 5998.799999999999
$ws.Range("M122").Value = -2060.5
$ws.Range("N122").Value = -10898.8
$ws.Range("H132").Value = 5561825
$ws.Range("I132").Value = 11769187
$ws.Range("J132").Value = 7869.5264
$ws.Range("K132").Value = 35307561
$ws.Range("L132").Value = 23608.5792
$ws.Range("M132").Value = -35305031
$ws.Range("N132").Value = -28668.5792
$ws.Range("H134").Value = 90196.2
$ws.Range("J134").Value = 90196.2
$ws.Range("L134").Value = 90196.2
$ws.Range("N134").Value = -100336.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 805.0909
$ws.Range("I99").Value = 805.6
$ws.Range("K99").Value = 805.6
$ws.Range("M99").Value = 692.4
$ws.Range("H105").Value = 500001900
$ws.Range("I105").Value = 500001900
$ws.Range("K105").Value = 500001900
$ws.Range("M105").Value = -500000153
$ws.Range("H134").Value = 7224.4546
$ws.Range("I134").Value = 5777.357
$ws.Range("K134").Value = 17332.071
$ws.Range("M134").Value = -14797.071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4366.6665
$ws.Range("I16").Value = 1550
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 1550
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = -1263
$ws.Range("N16").Value = -10574
$ws.Range("H22").Value = 515.6667
$ws.Range("I22").Value = 543.875
$ws.Range("J22").Value = 290
$ws.Range("K22").Value = 543.875
$ws.Range("L22").Value = 290
$ws.Range("M22").Value = -193.875
$ws.Range("N22").Value = -990
$ws.Range("H23").Value = 9966
$ws.Range("J23").Value = 9966
$ws.Range("L23").Value = 9966
$ws.Range("N23").Value = -10446
$ws.Range("H27").Value = 9966
$ws.Range("J27").Value = 9966
$ws.Range("L27").Value = 9966
$ws.Range("N27").Value = -10350
$ws.Range("H105").Value = 1652
$ws.Range("I105").Value = 1652
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1652
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 95
$ws.Range("H109").Value = 38994.54
$ws.Range("J109").Value = 38994.54
$ws.Range("L109").Value = 38994.54
$ws.Range("N109").Value = -41074.54
$ws.Range("H113").Value = 4366.6665
$ws.Range("I113").Value = 1550
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 1550
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = 620
$ws.Range("N113").Value = -14340
$ws.Range("H132").Value = 41385.188
$ws.Range("I132").Value = 4155.625
$ws.Range("K132").Value = 12466.875
$ws.Range("M132").Value = -9936.875
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 3208.4285
$ws.Range("I18").Value = 3391.8
$ws.Range("J18").Value = 2750
$ws.Range("K18").Value = 10175.4
$ws.Range("L18").Value = 8250
$ws.Range("M18").Value = -10006.4
$ws.Range("N18").Value = -8588
$ws.Range("H37").Value = 99947.14
$ws.Range("J37").Value = 99947.14
$ws.Range("L37").Value = 299841.42
$ws.Range("N37").Value = -300065.42
$ws.Range("H63").Value = 4860
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H66").Value = 4860
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H70").Value = 3059.75
$ws.Range("I70").Value = 3059.75
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 9179.25
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -8864.25
$ws.Range("H73").Value = 3059.75
$ws.Range("I73").Value = 3059.75
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 9179.25
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -8087.25
$ws.Range("H109").Value = 2294.7778
$ws.Range("I109").Value = 1831.625
$ws.Range("J109").Value = 6000
$ws.Range("K109").Value = 5494.875
$ws.Range("L109").Value = 18000
$ws.Range("M109").Value = -4454.875
$ws.Range("N109").Value = -20080
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2280
$ws.Range("I97").Value = 1950
$ws.Range("K97").Value = 1950
$ws.Range("M97").Value = -1454
$ws.Range("H122").Value = 4179.231
$ws.Range("I122").Value = 3513.7144
$ws.Range("K122").Value = 10541.1432
$ws.Range("M122").Value = -8091.143199999999
$ws.Range("H132").Value = 5960.8286
$ws.Range("I132").Value = 4153.091
$ws.Range("K132").Value = 12459.273
$ws.Range("M132").Value = -9929.273000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6764.423
$ws.Range("J46").Value = 7477.304
$ws.Range("L46").Value = 7477.304
$ws.Range("N46").Value = -7853.304
$ws.Range("H55").Value = 181.64
$ws.Range("I55").Value = 173.92308
$ws.Range("J55").Value = 190
$ws.Range("K55").Value = 173.92308
$ws.Range("L55").Value = 190
$ws.Range("M55").Value = -0.9230799999999988
$ws.Range("N55").Value = -536
$ws.Range("H93").Value = 10000
$ws.Range("J93").Value = 10000
$ws.Range("L93").Value = 10000
$ws.Range("N93").Value = -12496
$ws.Range("H94").Value = 23500
$ws.Range("J94").Value = 23500
$ws.Range("L94").Value = 23500
$ws.Range("N94").Value = -24852
$ws.Range("H100").Value = 31253318
$ws.Range("I100").Value = 125001850
$ws.Range("K100").Value = 125001850
$ws.Range("M100").Value = -125001309

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 33499.75
$ws.Range("J31").Value = 17999
$ws.Range("L31").Value = 17999
$ws.Range("N31").Value = -18695
$ws.Range("H96").Value = 2794
$ws.Range("I96").Value = 2739.1538
$ws.Range("J96").Value = 2865.3
$ws.Range("K96").Value = 2739.1538
$ws.Range("L96").Value = 2865.3
$ws.Range("M96").Value = -1366.1538
$ws.Range("N96").Value = -5611.3
$ws.Range("H100").Value = 1759.2858
$ws.Range("I100").Value = 1430.1818
$ws.Range("K100").Value = 2860.3636
$ws.Range("M100").Value = -2319.3636
$ws.Range("H109").Value = 51250
$ws.Range("J109").Value = 51250
$ws.Range("L109").Value = 51250
$ws.Range("N109").Value = -54024
$ws.Range("H126").Value = 6094.4443
$ws.Range("I126").Value = 6121.952
$ws.Range("K126").Value = 18365.856
$ws.Range("M126").Value = -15895.856
